$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 17: Inscritos 104 -> 106
$ws.Range("E17").Value = 106

# Row 19: Inscritos 50 -> 51, Pagos 25 -> 26, Inscricoes homologadas 25 -> 26
$ws.Range("E19").Value = 51
$ws.Range("F19").Value = 26
$ws.Range("H19").Value = 26

# Row 20: Inscritos 5 -> 6
$ws.Range("E20").Value = 6

# Row 24: Inscritos 20 -> 21, Pagos 12 -> 13, Inscricoes homologadas 12 -> 13
$ws.Range("E24").Value = 21
$ws.Range("F24").Value = 13
$ws.Range("H24").Value = 13

# Row 25: Inscritos 20 -> 21, Pagos 8 -> 9, Inscricoes homologadas 8 -> 9
$ws.Range("E25").Value = 21
$ws.Range("F25").Value = 9
$ws.Range("H25").Value = 9

# Row 34: Inscritos 17 -> 18
$ws.Range("E34").Value = 18

# Row 40: Inscritos 18 -> 19, Pagos 9 -> 10, Inscricoes homologadas 9 -> 10
$ws.Range("E40").Value = 19
$ws.Range("F40").Value = 10
$ws.Range("H40").Value = 10

# Row 42: Inscritos 31 -> 33
$ws.Range("E42").Value = 33

# Row 43: Inscritos 22 -> 24
$ws.Range("E43").Value = 24

# Row 68: Inscritos 13 -> 14
$ws.Range("E68").Value = 14

# Row 71: Pagos 13 -> 14, Inscricoes homologadas 13 -> 14
$ws.Range("F71").Value = 14
$ws.Range("H71").Value = 14

# Row 79: Inscritos 32 -> 34, Pagos 12 -> 13, Inscricoes homologadas 12 -> 13
$ws.Range("E79").Value = 34
$ws.Range("F79").Value = 13
$ws.Range("H79").Value = 13

$wb.Save()
